# Auto-generated Excel COM-interop script applying the diff changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2858.4211
$ws.Cells.Item(62, 9).Value = 2186
$ws.Cells.Item(62, 10).Value = 3605.5557
$ws.Cells.Item(62, 11).Value = 2186
$ws.Cells.Item(62, 12).Value = 3605.5557
$ws.Cells.Item(62, 13).Value = -1562
$ws.Cells.Item(62, 14).Value = -4853.5557

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 2858.4211
$ws.Cells.Item(65, 9).Value = 2186
$ws.Cells.Item(65, 10).Value = 3605.5557
$ws.Cells.Item(65, 11).Value = 10930
$ws.Cells.Item(65, 12).Value = 18027.7785
$ws.Cells.Item(65, 13).Value = -7810
$ws.Cells.Item(65, 14).Value = -24267.7785

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 5209999
$ws.Cells.Item(129, 10).Value = 1963.25
$ws.Cells.Item(129, 12).Value = 5889.75
$ws.Cells.Item(129, 14).Value = -15889.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 8005930.5
$ws.Cells.Item(132, 9).Value = 11117877
$ws.Cells.Item(132, 11).Value = 33353631
$ws.Cells.Item(132, 13).Value = -33351101

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(134, 8).Value = 29785.715
$ws.Cells.Item(134, 10).Value = 29785.715
$ws.Cells.Item(134, 12).Value = 29785.715
$ws.Cells.Item(134, 14).Value = -39925.715

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 4620.696
$ws.Cells.Item(138, 9).Value = 2639.8333
$ws.Cells.Item(138, 10).Value = 6781.636
$ws.Cells.Item(138, 11).Value = 7919.499899999999
$ws.Cells.Item(138, 12).Value = 20344.908
$ws.Cells.Item(138, 13).Value = -2779.499899999999
$ws.Cells.Item(138, 14).Value = -30624.908

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1827.2593
$ws.Cells.Item(45, 9).Value = 1084.8334
$ws.Cells.Item(45, 11).Value = 1084.8334
$ws.Cells.Item(45, 13).Value = -707.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 3607.2856
$ws.Cells.Item(122, 9).Value = 2250.2
$ws.Cells.Item(122, 11).Value = 6750.599999999999
$ws.Cells.Item(122, 13).Value = -4300.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(134, 8).Value = 35571.6
$ws.Cells.Item(134, 10).Value = 35571.6
$ws.Cells.Item(134, 12).Value = 35571.6
$ws.Cells.Item(134, 14).Value = -45711.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(57, 8).Value = 24465.25
$ws.Cells.Item(57, 10).Value = 24465.25
$ws.Cells.Item(57, 12).Value = 24465.25
$ws.Cells.Item(57, 14).Value = -25585.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 3905.1614
$ws.Cells.Item(132, 9).Value = 2306.25
$ws.Cells.Item(132, 10).Value = 5610.6665
$ws.Cells.Item(132, 11).Value = 6918.75
$ws.Cells.Item(132, 12).Value = 16831.9995
$ws.Cells.Item(132, 13).Value = -4388.75
$ws.Cells.Item(132, 14).Value = -21891.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1531.6923
$ws.Cells.Item(134, 9).Value = 919
$ws.Cells.Item(134, 10).Value = 2246.5
$ws.Cells.Item(134, 11).Value = 2757
$ws.Cells.Item(134, 12).Value = 6739.5
$ws.Cells.Item(134, 13).Value = -222
$ws.Cells.Item(134, 14).Value = -11809.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 2005.5385
$ws.Cells.Item(134, 9).Value = 1120.2222
$ws.Cells.Item(134, 10).Value = 3997.5
$ws.Cells.Item(134, 11).Value = 3360.6666
$ws.Cells.Item(134, 12).Value = 11992.5
$ws.Cells.Item(134, 13).Value = 1709.3334
$ws.Cells.Item(134, 14).Value = -22132.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2675.9167
$ws.Cells.Item(102, 9).Value = 1856.7778
$ws.Cells.Item(102, 13).Value = -234.7778000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 12401.6
$ws.Cells.Item(122, 9).Value = 20000
$ws.Cells.Item(122, 10).Value = 10502
$ws.Cells.Item(122, 11).Value = 60000
$ws.Cells.Item(122, 12).Value = 31506
$ws.Cells.Item(122, 13).Value = -57550
$ws.Cells.Item(122, 14).Value = -36406

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3021.3794
$ws.Cells.Item(132, 9).Value = 2534.5557
$ws.Cells.Item(132, 10).Value = 3818
$ws.Cells.Item(132, 11).Value = 7603.6671
$ws.Cells.Item(132, 12).Value = 11454
$ws.Cells.Item(132, 13).Value = -5073.6671
$ws.Cells.Item(132, 14).Value = -16514

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4568
$ws.Cells.Item(7, 9).Value = 1704
$ws.Cells.Item(7, 10).Value = 6000
$ws.Cells.Item(7, 11).Value = 1704
$ws.Cells.Item(7, 12).Value = 6000
$ws.Cells.Item(7, 13).Value = -1592
$ws.Cells.Item(7, 14).Value = -6224

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 111113160
$ws.Cells.Item(22, 9).Value = 250000720
$ws.Cells.Item(22, 10).Value = 3112.4
$ws.Cells.Item(22, 11).Value = 250000720
$ws.Cells.Item(22, 12).Value = 3112.4
$ws.Cells.Item(22, 13).Value = -250000425
$ws.Cells.Item(22, 14).Value = -3702.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 111113160
$ws.Cells.Item(27, 9).Value = 250000720
$ws.Cells.Item(27, 10).Value = 3112.4
$ws.Cells.Item(27, 11).Value = 250000720
$ws.Cells.Item(27, 12).Value = 3112.4
$ws.Cells.Item(27, 13).Value = -250000613
$ws.Cells.Item(27, 14).Value = -3326.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 2717.3794
$ws.Cells.Item(122, 9).Value = 2477.077
$ws.Cells.Item(122, 10).Value = 4800
$ws.Cells.Item(122, 11).Value = 7431.231000000001
$ws.Cells.Item(122, 12).Value = 14400
$ws.Cells.Item(122, 13).Value = -4981.231000000001
$ws.Cells.Item(122, 14).Value = -19300

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 4568
$ws.Cells.Item(126, 9).Value = 1704
$ws.Cells.Item(126, 10).Value = 6000
$ws.Cells.Item(126, 11).Value = 5112
$ws.Cells.Item(126, 12).Value = 18000
$ws.Cells.Item(126, 13).Value = -2642
$ws.Cells.Item(126, 14).Value = -22940

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(130, 8).Value = 25000
$ws.Cells.Item(130, 10).Value = 25000
$ws.Cells.Item(130, 12).Value = 25000
$ws.Cells.Item(130, 14).Value = -35040

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 668608.5600000001
$ws.Cells.Item(122, 9).Value = 715652
$ws.Cells.Item(122, 10).Value = 10000
$ws.Cells.Item(122, 11).Value = 2146956
$ws.Cells.Item(122, 12).Value = 30000
$ws.Cells.Item(122, 13).Value = -2144506
$ws.Cells.Item(122, 14).Value = -34900

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 8335276.5
$ws.Cells.Item(126, 9).Value = 1318.6666
$ws.Cells.Item(126, 10).Value = 16669234
$ws.Cells.Item(126, 11).Value = 3955.9998
$ws.Cells.Item(126, 12).Value = 50007702
$ws.Cells.Item(126, 13).Value = -1485.9998
$ws.Cells.Item(126, 14).Value = -50012642

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(131, 8).Value = 23000
$ws.Cells.Item(131, 10).Value = 23000
$ws.Cells.Item(131, 12).Value = 23000
$ws.Cells.Item(131, 14).Value = -33080

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 273063.16
$ws.Cells.Item(132, 9).Value = 348829.66
$ws.Cells.Item(132, 10).Value = 53340.3
$ws.Cells.Item(132, 11).Value = 1046488.98
$ws.Cells.Item(132, 12).Value = 160020.9
$ws.Cells.Item(132, 13).Value = -1043958.98
$ws.Cells.Item(132, 14).Value = -165080.9

